$d = $word.ActiveDocument

# --- Change 2: remove the _GoBack bookmark paragraph after "Von Johann Wolfgang von Martin"
# (that bookmark moves to the newly-authored content instead; the old spot becomes a plain
# empty paragraph).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Change 1: insert the new "Doenergedichte" section (two poems + a page-break + an
# empty formatted paragraph) at the very start of the document body.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="berschrift1"/>
      </w:pPr>
      <w:r>
        <w:t>Dönergedichte</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="berschrift2"/>
      </w:pPr>
      <w:r>
        <w:t>Martins Variante:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>An Laura :D</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Beiderseits regiert das Leid,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>rechts der Krieg und links der Neid.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Bieten Väterchen Frost die Stirn,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>bewiesen aber nicht zu viel Hirn.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Ein gutes doch das Ganze hat,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>macht die Wogen es doch wieder glatt.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Der Türkenfez für Freude steht,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>wenn manchmal er dein Messer wärmt.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Die Frauen dort recht offen,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>hat sie kaum getroffen,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Für Geld spüren sie deinen Drang,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>ist er auch manchmal ziemlich lang,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>manchmal breit,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>sie sind bereit.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Ein Meisterwerk aus Fleisch und Glück,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>keiner je es nicht erblickt,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Ziel so vieles Begehrens,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>kann es uns doch eines lehren.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Fällt das Atmen danach auch schwer,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>deine Unschuld erlangst du nimmermehr.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Dein erstes Mal vergisst du nimmer,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>es ist verankert tief im Innern.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Nie mehr kannst du davon lassen,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>musst du auch viel Geld verprassen,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>bereuen wirst du’s nie.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Hältst du es wieder in der Hand,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>rufst du es durchs ganze Land.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Das warme Fleisch, so oft gepriesen,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>die braune Haut, so sanft geblieben.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Leicht salzig schmeckt das Fleisch so zart,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>zu widerstehen ist mehr als hart.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Lässt du dich dann endlich gehen,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>kannst du es dann endlich sehen,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>die Öffnung voller Flüssigkeit,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>öffnest du sie lang und breit.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Haarlos glatt und wunderschön,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">lässt du ihn nie wieder </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>gehn</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>-Martin Schneider, heute irgendwann</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="berschrift2"/>
      </w:pPr>
      <w:r>
        <w:t>Michis Variante:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Zue</w:t>
      </w:r>
      <w:r>
        <w:t>rst kommt die Sauce feucht und w</w:t>
      </w:r>
      <w:r>
        <w:t>eiß,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>danach der Salat minderwertig aus Mais.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Dazu kommen Tomaten</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Gepflückt von Primaten!</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Reinrassiges Fleisch,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> braun wie der Kot.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">All das kommt zusammen im </w:t>
      </w:r>
      <w:r>
        <w:t>Fladenbrot</w:t>
      </w:r>
      <w:r>
        <w:t>!</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Welch</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> holdes </w:t>
      </w:r>
      <w:r>
        <w:t>Geschöpf</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> erstellet ihr da?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Sehet her und staunet</w:t>
      </w:r>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> der super </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Döna</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>!</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="berschrift1"/>
      </w:pPr>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="200"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="majorEastAsia" w:cstheme="majorBidi"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="36"/>
          <w:szCs w:val="28"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
    </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r = $d.Range(0, 0)
$r.InsertXML($xml)

Write-Host "Applied edits. Paragraphs now:" $d.Paragraphs.Count
